# Update the "how to use it ?" usage example for Ukarticles (G5) to match the
# revised UKNews API: the example now calls UKarticles.getNewsByCategory directly
# (with the swapped argument order) instead of going through getkeywordFromTitle.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUsage = '"in views.py":                                       from .UKNews import UKarticles                                      category=request.GET[''category'']
info_dict = UKarticles.getNewsByCategory(UKarticles,category)                                "in .html file":                                    {{ info_dict.title.0}}                          {{ info_dict.imageUrl.0}} '
$ws.Range("G5").Value = $newUsage

# Reflect the window/scroll state the author left the sheet in after the edit
# (best effort -- not all window-chrome properties round-trip through this host).
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 2
$ws.Range("G7").Select()
